$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record is published for "Puerro" (Vega Modelo de Temuco).
# It becomes the new row 8, pushing every existing record (old rows 8..155)
# down by one row; the oldest record (old row 155) ends up as the new row 156.
$ws.Rows.Item(8).Insert()

# Populate the newly inserted row 8 with the new weekly record.
$ws.Cells.Item(8, 1).Value = 10
$ws.Cells.Item(8, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(8, 3).Value = "La Araucanía"
$ws.Cells.Item(8, 4).Value = 44515
$ws.Cells.Item(8, 5).Value = 9
$ws.Cells.Item(8, 6).Value = 100112005
$ws.Cells.Item(8, 7).Value = "Puerro"
$ws.Cells.Item(8, 8).Value = "Azul de Maquehue"
$ws.Cells.Item(8, 9).Value = "Primera"
$ws.Cells.Item(8, 10).Value = 110
$ws.Cells.Item(8, 11).Value = 8000
$ws.Cells.Item(8, 12).Value = 8000
$ws.Cells.Item(8, 13).Value = 8000
$ws.Cells.Item(8, 14).Value = "`$/docena de paquetes"
$ws.Cells.Item(8, 15).Value = "Provincia de Cautín"
$ws.Cells.Item(8, 16).Value = 667
$ws.Cells.Item(8, 17).Value = 12
$ws.Cells.Item(8, 18).Value = "Hortaliza"
